$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# xlPasteAll = -4104, xlPasteFormats = -4122
$xlPasteFormats = -4122

# Insert two new blank rows at row 11 (row 10 "Mensajero" stays put; this
# reserves rows 11 and 12 for the new "Mensajero institución" / "Mensajero
# tutor" rows, while the old rows 11-13 shift down to 13-15).
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# ---------------------------------------------------------------------
# Row 10: "Mensajero alumno"
# ---------------------------------------------------------------------
# B10,C10,E10,F10,G10 held over "X" text from the old "Mensajero" row;
# clear that content out since the split-out row only keeps the D10 mark.
$ws.Range("B10:C10").ClearContents()
$ws.Range("E10:G10").ClearContents()

$ws.Cells.Item(10,1).Value = "Mensajero alumno"
$ws.Cells.Item(4,1).Copy()                 # donor A4 (style 9)
$ws.Cells.Item(10,1).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(4,7).Copy()                 # donor G4 (style 4)
$ws.Cells.Item(10,7).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(8,2).Copy()                 # donor B8 (style 7)
$ws.Cells.Item(10,2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10,3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10,5).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(10,6).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(10,4).Value = "X"
$ws.Cells.Item(3,2).Copy()                 # donor B3 (style 12)
$ws.Cells.Item(10,4).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(4,8).Copy()                 # donor H4 (style 32)
$ws.Cells.Item(10,8).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,9).Copy()                 # donor I4 (style 22)
$ws.Cells.Item(10,9).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,10).Copy()                # donor J4 (style 24)
$ws.Cells.Item(10,10).PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 11: "Mensajero institución"
# ---------------------------------------------------------------------
$ws.Cells.Item(11,1).Value = "Mensajero institución"
$ws.Cells.Item(4,1).Copy()                 # donor A4 (style 9)
$ws.Cells.Item(11,1).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(3,2).Copy()                 # donor B3 (style 12)
$ws.Cells.Item(11,2).Value = "X"
$ws.Cells.Item(11,2).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11,3).Value = "X"
$ws.Cells.Item(11,3).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11,6).Value = "X"
$ws.Cells.Item(11,6).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(8,2).Copy()                 # donor B8 (style 7)
$ws.Cells.Item(11,4).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(11,5).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(11,7).Value = "X"
$ws.Cells.Item(3,7).Copy()                 # donor G3 (style 15)
$ws.Cells.Item(11,7).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(4,8).Copy()                 # donor H4 (style 32)
$ws.Cells.Item(11,8).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,9).Copy()                 # donor I4 (style 22)
$ws.Cells.Item(11,9).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,10).Copy()                # donor J4 (style 24)
$ws.Cells.Item(11,10).PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 12: "Mensajero tutor"  (column C is intentionally left untouched /
# blank with no explicit format, matching the source workbook)
# ---------------------------------------------------------------------
$ws.Cells.Item(12,1).Value = "Mensajero tutor"
$ws.Cells.Item(4,1).Copy()                 # donor A4 (style 9)
$ws.Cells.Item(12,1).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(12,3).ClearFormats()

$ws.Cells.Item(8,2).Copy()                 # donor B8 (style 7)
$ws.Cells.Item(12,2).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(3,3).Copy()                 # donor C3 (style 1)
$ws.Cells.Item(12,4).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(12,6).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(12,5).Value = "X"
$ws.Cells.Item(3,2).Copy()                 # donor B3 (style 12)
$ws.Cells.Item(12,5).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(4,7).Copy()                 # donor G4 (style 4)
$ws.Cells.Item(12,7).PasteSpecial($xlPasteFormats)

$ws.Cells.Item(4,8).Copy()                 # donor H4 (style 32)
$ws.Cells.Item(12,8).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,9).Copy()                 # donor I4 (style 22)
$ws.Cells.Item(12,9).PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4,10).Copy()                # donor J4 (style 24)
$ws.Cells.Item(12,10).PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 13 (old row 11, "Diseñador encuestas", shifted down by the insert):
# only E13 picks up a different border style.
# ---------------------------------------------------------------------
$ws.Cells.Item(8,2).Copy()                 # donor B8 (style 7)
$ws.Cells.Item(13,5).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

$ws.Range("L2").Select()

Write-Output "done"
